$p = $ppt.ActivePresentation

# 1. Notes Master date placeholder: "22/01/2015" -> "01/07/15"
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$hf.DateAndTime.Text = "01/07/15"

# 2. Slide 1 - remove the "January 2015" paragraph from the content placeholder,
#    leaving a single empty paragraph.
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$tr1.Paragraphs(1, 1).Delete()

# 3. Slide 7 - merge "Measurements " + "(top bar)" runs into a single run.
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(1)
$tr7 = $shp7.TextFrame.TextRange
$run7a = $tr7.Characters(1, 13)
$run7a.Text = "Measurements (top bar)"
$tr7.Characters(23, 9).Delete()

# 4. Slide 9 - split "#ff08518" into "#" + "f08518" (fixing the typo'd hex value).
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(8)
$tr9 = $shp9.TextFrame.TextRange
$tr9.Characters(2, 7).Text = "f08518"
